$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.363.38"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = "'2.622.75"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'596.99"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = "'152.77"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = "'0.553"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'2.619.28"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('E10').Value = '  -2.05%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').Value = "'27.60"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = "'3.093.78"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('E16').Value = '  -2.41%  '
$ws.Range('D17').Value = "'67.413.11"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = "'2.616.18"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'11.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('D20').Value = "'362.87"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.87%  '
$ws.Range('D21').Value = "'7.48"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.03%  '
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('E23').Value = '  +3.98%  '
$ws.Range('D24').Value = "'1.00"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = "'70.99"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.84%  '
$ws.Range('D26').Value = "'10.09"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.58%  '
$ws.Range('D27').Value = "'2.761.13"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').Value = "'584.28"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.13%  '
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').Value = "'1.40"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.99%  '
$ws.Range('D32').Value = "'7.82"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.16%  '
$ws.Range('E33').Value = '  -0.98%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  -5.97%  '
$ws.Range('E36').Value = '  -2.09%  '
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('D38').Value = "'157.22"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.67%  '
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = "'5.27"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('D44').Value = "'41.17"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = "'16.34"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('D47').Value = "'156.45"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').Value = "'0.0₆0287"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').Value = "'20.62"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.63%  '
